# ActionPlan.xlsx update — "check, checkmate and stalemate"
#
# The action-plan sheet tracks one commit per row in column A (bold text),
# with an "x" in column B once that commit/task is done. This change:
#   1) marks the previously-last task (row 208, "Implement check, checkmate
#      and stalemate") as done by putting an "x" in B208, and
#   2) appends the next batch of completed commit-message rows (every other
#      row, starting at row 210) for the new work items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Mark row 208 as done.
$ws.Cells.Item(208, 2).Value = "x"

# 2) Append the new task rows, two sheet-rows apart, starting at row 210.
$newEntries = @(
    "Introduce Situation class",
    "Test check, checkmate and stalemate",
    "Add Situation playMoves functions",
    "Test complete games",
    "Add castle tests",
    "Reorganizing moving rules",
    "Rename trajectory to longRange",
    "Test and implement standard chess kingside castle",
    "More castle tests",
    "Attempt to generalize castle side",
    "Implement and test conditioned pos vectors",
    "Remove unused pos functions"
)

$row = 210
foreach ($entry in $newEntries) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $entry
    $cell.Font.Bold = $true
    $row += 2
}

# "Add Situation playMoves functions" (row 214) is a section-style header:
# bold + wrap text, matching the other section headers already in the
# sheet (e.g. A66, A87).
$ws.Cells.Item(214, 1).WrapText = $true

# Move the view/selection down to the newly added rows, like the author did.
$ws.Application.ActiveWindow.ScrollRow = 200
$ws.Cells.Item(210, 2).Select()
